$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 137.5
$ws.Range("I12").Value = 133.33333
$ws.Range("K12").Value = 133.33333
$ws.Range("M12").Value = 36.66667000000001
$ws.Range("H32").Value = 2493.5
$ws.Range("I32").Value = 1793.4
$ws.Range("K32").Value = 1793.4
$ws.Range("M32").Value = -1467.4
$ws.Range("H103").Value = 624.0625
$ws.Range("J103").Value = 641.7
$ws.Range("L103").Value = 1925.1
$ws.Range("N103").Value = -3097.1
$ws.Range("H112").Value = 2224.75
$ws.Range("J112").Value = 2224.75
$ws.Range("L112").Value = 6674.25
$ws.Range("N112").Value = -8890.25
$ws.Range("H123").Value = 100000
$ws.Range("J123").Value = 100000
$ws.Range("L123").Value = 100000
$ws.Range("N123").Value = -109800
$ws.Range("H138").Value = 2504.29
$ws.Range("I138").Value = 2128.4211
$ws.Range("J138").Value = 2734.6614
$ws.Range("K138").Value = 6385.263300000001
$ws.Range("L138").Value = 8203.984199999999
$ws.Range("M138").Value = -1245.263300000001
$ws.Range("N138").Value = -18483.9842

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 779.43475
$ws.Range("I2").Value = 663.36365
$ws.Range("J2").Value = 3333
$ws.Range("K2").Value = 663.36365
$ws.Range("L2").Value = 3333
$ws.Range("M2").Value = -550.36365
$ws.Range("N2").Value = -3559
$ws.Range("H45").Value = 1656.5714
$ws.Range("I45").Value = 1639.4
$ws.Range("K45").Value = 1639.4
$ws.Range("M45").Value = -1262.4
$ws.Range("H97").Value = 2154.5715
$ws.Range("I97").Value = 847
$ws.Range("K97").Value = 847
$ws.Range("M97").Value = -351
$ws.Range("H116").Value = 779.43475
$ws.Range("I116").Value = 663.36365
$ws.Range("J116").Value = 3333
$ws.Range("K116").Value = 663.36365
$ws.Range("L116").Value = 3333
$ws.Range("M116").Value = 1630.63635
$ws.Range("N116").Value = -7921
$ws.Range("H130").Value = 18666.334
$ws.Range("J130").Value = 18666.334
$ws.Range("L130").Value = 18666.334
$ws.Range("N130").Value = -28706.334
$ws.Range("H132").Value = 2458.52
$ws.Range("I132").Value = 2275.5908
$ws.Range("K132").Value = 6826.7724
$ws.Range("M132").Value = -4296.7724

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 779.43475
$ws.Range("I3").Value = 663.36365
$ws.Range("J3").Value = 3333
$ws.Range("K3").Value = 663.36365
$ws.Range("L3").Value = 3333
$ws.Range("M3").Value = -549.36365
$ws.Range("N3").Value = -3561
$ws.Range("H20").Value = 7571.2856
$ws.Range("I20").Value = 3831.6667
$ws.Range("K20").Value = 3831.6667
$ws.Range("M20").Value = -3584.6667
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("N49").Value = 0
$ws.Range("H105").Value = 2949.8823
$ws.Range("I105").Value = 2712.8333
$ws.Range("J105").Value = 3518.8
$ws.Range("K105").Value = 2712.8333
$ws.Range("L105").Value = 3518.8
$ws.Range("M105").Value = -965.8332999999998
$ws.Range("N105").Value = -7012.8
$ws.Range("H134").Value = 1738.1428
$ws.Range("I134").Value = 1394.5555
$ws.Range("J134").Value = 3799.6667
$ws.Range("K134").Value = 4183.666499999999
$ws.Range("L134").Value = 11399.0001
$ws.Range("M134").Value = -1648.666499999999
$ws.Range("N134").Value = -16469.0001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66826.2
$ws.Range("I7").Value = 91053.82000000001
$ws.Range("K7").Value = 91053.82000000001
$ws.Range("M7").Value = -90940.82000000001
$ws.Range("H94").Value = 6941.3335
$ws.Range("J94").Value = 3000
$ws.Range("L94").Value = 3000
$ws.Range("N94").Value = -3902

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 545.6667
$ws.Range("I7").Value = 545.6667
$ws.Range("K7").Value = 1637.0001
$ws.Range("M7").Value = -1525.0001
$ws.Range("H8").Value = 329.42856
$ws.Range("I8").Value = 329.42856
$ws.Range("K8").Value = 988.28568
$ws.Range("M8").Value = -849.28568
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("H92").Value = 285.3
$ws.Range("J92").Value = 285.44446
$ws.Range("L92").Value = 856.33338
$ws.Range("N92").Value = -3352.33338
$ws.Range("H131").Value = 32038.055
$ws.Range("I131").Value = 139803.25
$ws.Range("J131").Value = 2309.724
$ws.Range("K131").Value = 419409.75
$ws.Range("L131").Value = 6929.172
$ws.Range("M131").Value = -414369.75
$ws.Range("N131").Value = -17009.172

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("H40").Value = 50000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H63").Value = 32557
$ws.Range("I63").Value = 15000
$ws.Range("K63").Value = 15000
$ws.Range("M63").Value = -14314
$ws.Range("H66").Value = 32557
$ws.Range("I66").Value = 15000
$ws.Range("K66").Value = 45000
$ws.Range("M66").Value = -41568
$ws.Range("H70").Value = 8476.666999999999
$ws.Range("J70").Value = 8973.5
$ws.Range("L70").Value = 8973.5
$ws.Range("N70").Value = -9513.5
$ws.Range("H73").Value = 8476.666999999999
$ws.Range("J73").Value = 8973.5
$ws.Range("L73").Value = 8973.5
$ws.Range("N73").Value = -10845.5
$ws.Range("H80").Value = 17966
$ws.Range("I80").Value = 4811.6
$ws.Range("J80").Value = 26187.5
$ws.Range("K80").Value = 4811.6
$ws.Range("L80").Value = 26187.5
$ws.Range("M80").Value = -3813.6
$ws.Range("N80").Value = -28183.5
$ws.Range("H83").Value = 17966
$ws.Range("I83").Value = 4811.6
$ws.Range("J83").Value = 26187.5
$ws.Range("K83").Value = 24058
$ws.Range("L83").Value = 130937.5
$ws.Range("M83").Value = -19066
$ws.Range("N83").Value = -140921.5
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 300000
$ws.Range("N134").Value = -305070
$ws.Range("H136").Value = 54999
$ws.Range("J136").Value = 54999
$ws.Range("L136").Value = 164997
$ws.Range("N136").Value = -170097

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 994
$ws.Range("I22").Value = 994
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 994
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = -699
$ws.Range("H27").Value = 994
$ws.Range("I27").Value = 994
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 994
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = -887
$ws.Range("H33").Value = 51000
$ws.Range("J33").Value = 51000
$ws.Range("L33").Value = 51000
$ws.Range("N33").Value = -51580
$ws.Range("H43").Value = 6399007
$ws.Range("J43").Value = 6937086.5
$ws.Range("L43").Value = 6937086.5
$ws.Range("N43").Value = -6937472.5
$ws.Range("H46").Value = 2210.625
$ws.Range("I46").Value = 1409.4286
$ws.Range("J46").Value = 2833.7778
$ws.Range("K46").Value = 1409.4286
$ws.Range("L46").Value = 2833.7778
$ws.Range("M46").Value = -1221.4286
$ws.Range("N46").Value = -3209.7778
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("N113").Value = 0
$ws.Range("H122").Value = 4591.0625
$ws.Range("J122").Value = 2852.5
$ws.Range("L122").Value = 8557.5
$ws.Range("N122").Value = -13457.5
$ws.Range("H136").Value = 3006.3333
$ws.Range("I136").Value = 3316.4614
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 9949.3842
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -7399.3842
$ws.Range("N136").Value = -11700

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 14999.667
$ws.Range("I29").Value = 10799.2
$ws.Range("J29").Value = 18000
$ws.Range("K29").Value = 10799.2
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = -10509.2
$ws.Range("N29").Value = -18580
$ws.Range("H107").Value = 1215.7222
$ws.Range("J107").Value = 1278.1428
$ws.Range("L107").Value = 3834.4284
$ws.Range("N107").Value = -7674.428400000001
$ws.Range("H126").Value = 2201.5
$ws.Range("I126").Value = 2294.8
$ws.Range("K126").Value = 6884.400000000001
$ws.Range("M126").Value = -4414.400000000001
$ws.Range("H132").Value = 1269.579
$ws.Range("I132").Value = 1157.9375
$ws.Range("J132").Value = 1865
$ws.Range("K132").Value = 3473.8125
$ws.Range("L132").Value = 5595
$ws.Range("M132").Value = -943.8125
$ws.Range("N132").Value = -10655
$ws.Range("H136").Value = 1215.6818
$ws.Range("J136").Value = 1844.3334
$ws.Range("L136").Value = 5533.0002
$ws.Range("N136").Value = -10633.0002
